$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 46031
$ws.Range("B2").Value = 13699.5217630316
$ws.Range("C2").Value = 12158.2520404986
$ws.Range("D2").Value = 15819.86
$ws.Range("E2").Value = 9494.3924942691
$ws.Range("F2").Value = 243.032688948654

$ws.Range("A3").Value = 46032
$ws.Range("B3").Value = 5457.36453438434
$ws.Range("C3").Value = 8586.87244125754
$ws.Range("D3").Value = 12075.86
$ws.Range("E3").Value = 8792.06432695207
$ws.Range("F3").Value = 220.961532008734

$ws.Range("A4").Value = 46033
$ws.Range("B4").Value = 5322.78549778354
$ws.Range("C4").Value = 8865.98964384413
$ws.Range("D4").Value = 12075.86
$ws.Range("E4").Value = 8773.66341800627
$ws.Range("F4").Value = 231.824710910433

$ws.Range("A5").Value = 46034
$ws.Range("B5").Value = 13347.132397846
$ws.Range("C5").Value = 13341.5778370522
$ws.Range("D5").Value = 12075.86
$ws.Range("E5").Value = 9161.46092255359
$ws.Range("F5").Value = 434.46578165024

$ws.Range("A6").Value = 46035
$ws.Range("B6").Value = 12986.4338904735
$ws.Range("C6").Value = 12772.8233416059
$ws.Range("D6").Value = 12075.86
$ws.Range("E6").Value = 8869.93228138575
$ws.Range("F6").Value = 398.620650957987

$ws.Range("A7").Value = 46036
$ws.Range("B7").Value = 13015.7485041368
$ws.Range("C7").Value = 12219.836935601
$ws.Range("D7").Value = 12075.86
$ws.Range("E7").Value = 8868.63801881166
$ws.Range("F7").Value = 375.525623100529

$ws.Range("A8").Value = 46037
$ws.Range("B8").Value = 13015.7485041368
$ws.Range("C8").Value = 12107.8336665116
$ws.Range("D8").Value = 12075.86
$ws.Range("E8").Value = 8868.63801881166
$ws.Range("F8").Value = 370.858820221801

$ws.Range("A9").Value = 46038
$ws.Range("B9").Value = 13015.7485041368
$ws.Range("C9").Value = 11380.3994452514
$ws.Range("D9").Value = 12075.86
$ws.Range("E9").Value = 8868.59983477215
$ws.Range("F9").Value = 340.547470000983

$ws.Range("A10").Value = 46039
$ws.Range("B10").Value = 5247.16359875529
$ws.Range("C10").Value = 7686.46507063572
$ws.Range("D10").Value = 12075.86
$ws.Range("E10").Value = 8468.05208663344
$ws.Range("F10").Value = 169.944048219548

$ws.Range("A11").Value = 46040
$ws.Range("B11").Value = 5137.84212801753
$ws.Range("C11").Value = 7675.0370428071
$ws.Range("D11").Value = 12075.86
$ws.Range("E11").Value = 8460.20612956141
$ws.Range("F11").Value = 169.140965515354

$ws.Range("A12").Value = 46041
$ws.Range("B12").Value = 12719.4875548453
$ws.Range("C12").Value = 11557.6143661175
$ws.Range("D12").Value = 12075.86
$ws.Range("E12").Value = 8648.2987964125
$ws.Range("F12").Value = 338.752215105418

$ws.Range("A13").Value = 46042
$ws.Range("B13").Value = 12719.4875548453
$ws.Range("C13").Value = 11753.6153794981
$ws.Range("D13").Value = 12075.86
$ws.Range("E13").Value = 8648.2987964125
$ws.Range("F13").Value = 346.918923996275

$ws.Range("A14").Value = 46043
$ws.Range("B14").Value = 12719.4875548453
$ws.Range("C14").Value = 11918.1564567915
$ws.Range("D14").Value = 12075.86
$ws.Range("E14").Value = 8648.2987964125
$ws.Range("F14").Value = 353.774802216835

$ws.Range("A15").Value = 46044
$ws.Range("B15").Value = 12719.4875548453
$ws.Range("C15").Value = 12117.8552398677
$ws.Range("D15").Value = 12075.86
$ws.Range("E15").Value = 8648.2987964125
$ws.Range("F15").Value = 362.095584845008

